# Generate Report for Handback
# Update the "xliff generate" / handback timestamps that were refreshed
# when the report was regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for a5f17f62-... row
$wsOverview.Range("G4").Value = "2016-08-28 10:53:58"

# zh-cn sheet: Correspond Handoff / Handback Datetime for a5f17f62-... row
$wsZhCn.Range("H4").Value = "2016-08-28 10:53:54"
$wsZhCn.Range("K4").Value = "2016-08-28 10:54:16"

# de-de sheet: Correspond Handback Datetime for a5f17f62-... row
$wsDeDe.Range("K4").Value = "2016-08-28 10:54:22"
